$d = $word.ActiveDocument

$pkgPrefix = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgSuffix = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# -------------------------------------------------------------------------
# 1) "logo.nologo vt.global_cursor_default=0" paragraph: drop the proofErr
#    wrappers, split the runs differently and move the _GoBack bookmark
#    here (right after "logo.nologo").
# -------------------------------------------------------------------------
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "logo.nologo*global_cursor_default*") {
        $target1 = $p
        break
    }
}

$body1 = "<w:body><w:p>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>logo.nologo</w:t></w:r>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
  "<w:bookmarkEnd w:id=`"0`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`"> vt.globa</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>l</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>_cursor_default=0</w:t></w:r>" +
  "</w:p></w:body>"

$null = $target1.Range.InsertXML($pkgPrefix + $body1 + $pkgSuffix)

# -------------------------------------------------------------------------
# 2) "wifi-country.service" paragraph: remove the _GoBack bookmark that
#    used to live here (it moved to the paragraph above).
# -------------------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*wifi-country.service*") {
        $target2 = $p
        break
    }
}

$body2 = "<w:body><w:p>" +
  "<w:r><w:t>&#8220;</w:t></w:r>" +
  "<w:r><w:tab/></w:r>" +
  "<w:r><w:tab/></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:tab/><w:t xml:space=`"preserve`">&#8220; </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>wifi</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:t>-country.service</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "</w:p></w:body>"

$null = $target2.Range.InsertXML($pkgPrefix + $body2 + $pkgSuffix)

# -------------------------------------------------------------------------
# 3) Add the new "hciuart.service" paragraph right after the
#    "fake-hwclock.service" paragraph (replacing the first of the blank
#    paragraphs that follows it).
# -------------------------------------------------------------------------
$target3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*fake-hwclock.service*") {
        $target3 = $p.Next()
        break
    }
}

$body3 = "<w:body><w:p>" +
  "<w:r><w:t>&#8220;</w:t></w:r>" +
  "<w:r><w:tab/></w:r>" +
  "<w:r><w:tab/></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:tab/></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`">&#8220; </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>hciuart</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>.service</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>&#8221;</w:t></w:r>" +
  "</w:p></w:body>"

$null = $target3.Range.InsertXML($pkgPrefix + $body3 + $pkgSuffix)

Write-Output "Applied 3 edits: logo.nologo split/bookmark move, wifi-country bookmark removal, hciuart.service paragraph insert."
